$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: G3, H3 -> 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: D4, E4 -> 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: H5 -> 1
$ws.Range("H5").Value = 1

# Row 6: H6 -> 1
$ws.Range("H6").Value = 1

# Row 7: H7 -> 1
$ws.Range("H7").Value = 1

# Row 8: H8 -> 1
$ws.Range("H8").Value = 1

# Row 9: H9 -> 1
$ws.Range("H9").Value = 1

# Row 10: H10 -> 1
$ws.Range("H10").Value = 1

# Row 11: D11, E11 -> 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1

# Row 12: H12 -> 1
$ws.Range("H12").Value = 1

# Row 13: D13, E13 -> 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1

# Row 14: H14 -> 1
$ws.Range("H14").Value = 1

# Row 15: H15 -> 1
$ws.Range("H15").Value = 1

# Row 16: H16 -> 1
$ws.Range("H16").Value = 1

# Row 17: H17 -> 1
$ws.Range("H17").Value = 1

# Row 18: H18 -> 1
$ws.Range("H18").Value = 1
